# Initial Data File Update
# Adds 5 new transaction rows (152-156) to the "Transacciones" sheet,
# mirroring new grocery / household purchases recorded for 06-07 May 2019.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transacciones")

# ---------------------------------------------------------------------------
# Row 152
# ---------------------------------------------------------------------------
$ws.Range("A151").Copy()
$ws.Range("A152").PasteSpecial(-4122)
$ws.Range("A152").Value = 43591
$ws.Range("B152").Value = 17.5
$ws.Range("C152").Value = "Bolsa para basura"
$ws.Range("D152").Value = "Limpieza"
$ws.Range("E152").Value = "Gasto"
$ws.Range("F152").Value = "Tarjeta banamex"
$ws.Range("G152").Value = "Soriana"
$ws.Range("K152").Formula = "=K151-B152"
$ws.Range("L152").Value = 3547.12
$ws.Range("M152").Value = 53
$ws.Range("N152").Formula = "=SUM(K152:M152)"
$ws.Range("O152").Formula = "=N152-4000"
$ws.Range("P152").Style = "Bueno"
$ws.Range("P152").Formula = '=O152-Ahorros!$E$4'

# ---------------------------------------------------------------------------
# Row 153
# ---------------------------------------------------------------------------
$ws.Range("A151").Copy()
$ws.Range("A153").PasteSpecial(-4122)
$ws.Range("A153").Value = 43591
$ws.Range("B153").Value = 15.6
$ws.Range("C153").Value = "Chips Jalapeño"
$ws.Range("D153").Value = "Golosinas"
$ws.Range("E153").Value = "Gasto"
$ws.Range("F153").Value = "Tarjeta banamex"
$ws.Range("G153").Value = "Soriana"
$ws.Range("K153").Formula = "=K152-B153"
$ws.Range("L153").Value = 3547.12
$ws.Range("M153").Value = 53
$ws.Range("N153").Formula = "=SUM(K153:M153)"
$ws.Range("O153").Formula = "=N153-4000"
$ws.Range("P153").Style = "Bueno"
$ws.Range("P153").Formula = '=O153-Ahorros!$E$4'

# ---------------------------------------------------------------------------
# Row 154
# ---------------------------------------------------------------------------
$ws.Range("A151").Copy()
$ws.Range("A154").PasteSpecial(-4122)
$ws.Range("A154").Value = 43591
$ws.Range("B154").Value = 31.3
$ws.Range("C154").Value = "Papaya"
$ws.Range("D154").Value = "Despensa"
$ws.Range("E154").Value = "Gasto"
$ws.Range("F154").Value = "Tarjeta Banamex"
$ws.Range("G154").Value = "Soriana"
$ws.Range("K154").Formula = "=K153-B154"
$ws.Range("L154").Value = 3547.12
$ws.Range("M154").Value = 53
$ws.Range("N154").Formula = "=SUM(K154:M154)"
$ws.Range("O154").Formula = "=N154-4000"
$ws.Range("P154").Style = "Bueno"
$ws.Range("P154").Formula = '=O154-Ahorros!$E$4'

# ---------------------------------------------------------------------------
# Row 155
# ---------------------------------------------------------------------------
$ws.Range("A151").Copy()
$ws.Range("A155").PasteSpecial(-4122)
$ws.Range("A155").Value = 43592
$ws.Range("B155").Value = 9
$ws.Range("C155").Value = "Dr. Pepper"
$ws.Range("D155").Value = "Despensa"
$ws.Range("E155").Value = "Gasto"
$ws.Range("F155").Value = "Tarjeta Santander"
$ws.Range("G155").Value = "Extra"
$ws.Range("K155").Value = 5064.18
$ws.Range("L155").Formula = "=L154-B155"
$ws.Range("M155").Value = 53
$ws.Range("N155").Formula = "=SUM(K155:M155)"
$ws.Range("O155").Formula = "=N155-4000"
$ws.Range("P155").Style = "Bueno"
$ws.Range("P155").Formula = '=O155-Ahorros!$E$4'

# ---------------------------------------------------------------------------
# Row 156
# ---------------------------------------------------------------------------
$ws.Range("A151").Copy()
$ws.Range("A156").PasteSpecial(-4122)
$ws.Range("A156").Value = 43592
$ws.Range("B156").Value = 79
$ws.Range("C156").Value = "Lego Personaje"
$ws.Range("D156").Value = "Juguetes"
$ws.Range("E156").Value = "Gasto"
$ws.Range("F156").Value = "Tarjeta Santander"
$ws.Range("G156").Value = "Lego Store"
$ws.Range("K156").Value = 5064.18
$ws.Range("L156").Formula = "=L155-B156"
$ws.Range("M156").Value = 53
$ws.Range("N156").Formula = "=SUM(K156:M156)"
$ws.Range("O156").Formula = "=N156-4000"
$ws.Range("P156").Style = "Bueno"
$ws.Range("P156").Formula = '=O156-Ahorros!$E$4'

# ---------------------------------------------------------------------------
# Update the saved selection/view state to reflect the new last row
# ---------------------------------------------------------------------------
$ws.Range("A157").Select()
